# Updated symbol list on Fri Dec 30 08:58:22 UTC 2022 with GitHub Actions
#
# All "Price" (column D) values are stored as TEXT in this workbook (not
# numbers), so numeric-looking strings are written with a leading apostrophe
# to force text entry (mirrors how real Excel avoids auto-converting them to
# numbers), then the cell style is reset back to "Normal" so no stray
# number-format/quote-prefix style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice([string]$addr, [string]$val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - BNB
Set-TextPrice "D2" "244.02"

# Row 3 - OKB
Set-TextPrice "D3" "23.82"

# Row 4 - HuobiToken
Set-TextPrice "D4" "5.133"

# Row 5 - Cronos
Set-TextPrice "D5" "0.05748"

# Row 6 - KuCoinToken
Set-TextPrice "D6" "6.479"

# Row 7 - GateToken
Set-TextPrice "D7" "3.138"

# Row 9 - FTXToken
Set-TextPrice "D9" "0.8387"

# Row 10 - WazirX
Set-TextPrice "D10" "0.1339"

# Row 11 - MandalaExchangeToken
Set-TextPrice "D11" "0.06940"

# Row 12 - LiechtensteinCryptoassetsExchange
Set-TextPrice "D12" "0.03130"

# Row 13 - BitrueCoin
Set-TextPrice "D13" "0.02847"

# Row 14 - BitMartToken
Set-TextPrice "D14" "0.09364"

# Row 15 - MCDex
Set-TextPrice "D15" "3.739"

# Row 16 - BitForexToken
Set-TextPrice "D16" "0.001514"

# Row 17 - CoinExToken
Set-TextPrice "D17" "0.04652"

# Row 18 - One (price + it becomes the new "Worst in 24h")
Set-TextPrice "D18" "0.0005990"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# Row 19 - TigerCash
Set-TextPrice "D19" "0.006075"

# Row 20 - BitKan
Set-TextPrice "D20" "0.001234"

# Row 21 - HotbitToken
Set-TextPrice "D21" "0.004276"

# Row 22 - NitroEx
Set-TextPrice "D22" "0.00008697"

# Row 23 - LEO
Set-TextPrice "D23" "3.502"

# Row 25 - BitpandaEcosystemToken
Set-TextPrice "D25" "0.3174"

# Row 40 - IDEX
Set-TextPrice "D40" "0.03613"

# Rows 41-43 - coins reshuffled (BKEXToken/CEJI/KickToken rotate down one slot)
# Row 41 becomes KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextPrice "D41" "0.006378"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 becomes BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextPrice "D42" "0.1049"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 becomes CEJI (keeps the "Best in 24h" marker)
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextPrice "D43" "0.003299"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"

# Row 44 - LocalTraders
Set-TextPrice "D44" "0.007378"

# Row 48 - BOLO
Set-TextPrice "D48" "0.002278"

# Row 49 - CryptobidCoin
Set-TextPrice "D49" "0.00002099"
